$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns (U = Paineis, V = Ferramentas) ---
# Copy the header style from the last existing header cell (T1) so the
# new header cells keep the same bold/centered/bordered formatting.
$ws.Range("T1").Copy()
$ws.Range("U1:V1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("U1").Value = "Painéis"
$ws.Range("V1").Value = "Ferramentas"

# --- New row 14 (new survey submission) ---
$ws.Range("A14").Value = "jessica.mariano@mrv.com.br"
$ws.Range("I14").Value = "2025-05-20 16:24:13"
$ws.Range("U14").Value = "Painel do Portifólio - Planejamento da Produção - PLNESROBR004: Comentário Painel 1`n; Painel Produção Produtividade e MO - PLNESROBR005: Comentário Painel 2"
$ws.Range("V14").Value = "Controle de concretagem :Controlar quantidade e dias de concretagens da obra,Excel,💎 Muito Importante,6.0; Planilha de medição de EMP:Medição de Empreiteiros,Excel,🪙 Importante,4.0"

# The multi-line text entered above makes the engine auto-expand the row
# height; restore the default (non-custom) row height to match a plain
# data write.
$ws.Rows.Item(14).AutoFit()
